$wb = $excel.ActiveWorkbook

# --- Sheet "day": convert D797:D801 bsecode values from text to numeric ---
$daySheet = $wb.Worksheets.Item("day")
$daySheet.Cells.Item(797, 4).Value = 500027
$daySheet.Cells.Item(798, 4).Value = 539524
$daySheet.Cells.Item(799, 4).Value = 512599
$daySheet.Cells.Item(800, 4).Value = 542650
$daySheet.Cells.Item(801, 4).Value = 500257

# --- Sheet "week": append new rows 414-436 ---
$weekSheet = $wb.Worksheets.Item("week")

# Row 414
$weekSheet.Cells.Item(414, 1).Value = 1
$weekSheet.Cells.Item(414, 2).Value = "OFSS"
$weekSheet.Cells.Item(414, 3).Value = "Oracle Financial Services Software Limited"
$weekSheet.Cells.Item(414, 4).Value = 532466
$weekSheet.Cells.Item(414, 5).Value = 1.93
$weekSheet.Cells.Item(414, 6).Value = 10888.05
$weekSheet.Cells.Item(414, 7).Value = 432097
$weekSheet.Cells.Item(414, 8).Value = "week"
$weekSheet.Cells.Item(414, 9).Value = "25/10/2024 11:35:48"

# Row 415
$weekSheet.Cells.Item(415, 1).Value = 2
$weekSheet.Cells.Item(415, 2).Value = "ABB"
$weekSheet.Cells.Item(415, 3).Value = "Abb India Limited"
$weekSheet.Cells.Item(415, 4).Value = 500002
$weekSheet.Cells.Item(415, 5).Value = -1.65
$weekSheet.Cells.Item(415, 6).Value = 7521.65
$weekSheet.Cells.Item(415, 7).Value = 287323
$weekSheet.Cells.Item(415, 8).Value = "week"
$weekSheet.Cells.Item(415, 9).Value = "25/10/2024 11:35:48"

# Row 416
$weekSheet.Cells.Item(416, 1).Value = 3
$weekSheet.Cells.Item(416, 2).Value = "LTIM"
$weekSheet.Cells.Item(416, 3).Value = "LTI Mindtree Ltd"
$weekSheet.Cells.Item(416, 4).Value = 540005
$weekSheet.Cells.Item(416, 5).Value = -1.12
$weekSheet.Cells.Item(416, 6).Value = 5903.2
$weekSheet.Cells.Item(416, 7).Value = 191734
$weekSheet.Cells.Item(416, 8).Value = "week"
$weekSheet.Cells.Item(416, 9).Value = "25/10/2024 11:35:48"

# Row 417
$weekSheet.Cells.Item(417, 1).Value = 4
$weekSheet.Cells.Item(417, 2).Value = "PIIND"
$weekSheet.Cells.Item(417, 3).Value = "Pi Industries Limited"
$weekSheet.Cells.Item(417, 4).Value = 523642
$weekSheet.Cells.Item(417, 5).Value = 0.16
$weekSheet.Cells.Item(417, 6).Value = 4327.65
$weekSheet.Cells.Item(417, 7).Value = 74396
$weekSheet.Cells.Item(417, 8).Value = "week"
$weekSheet.Cells.Item(417, 9).Value = "25/10/2024 11:35:48"

# Row 418
$weekSheet.Cells.Item(418, 1).Value = 5
$weekSheet.Cells.Item(418, 2).Value = "NAVINFLUOR"
$weekSheet.Cells.Item(418, 3).Value = "Navin Fluorine International Limited"
$weekSheet.Cells.Item(418, 4).Value = 532504
$weekSheet.Cells.Item(418, 5).Value = -3.22
$weekSheet.Cells.Item(418, 6).Value = 3301.75
$weekSheet.Cells.Item(418, 7).Value = 234141
$weekSheet.Cells.Item(418, 8).Value = "week"
$weekSheet.Cells.Item(418, 9).Value = "25/10/2024 11:35:48"

# Row 419
$weekSheet.Cells.Item(419, 1).Value = 6
$weekSheet.Cells.Item(419, 2).Value = "MPHASIS"
$weekSheet.Cells.Item(419, 3).Value = "Mphasis Limited"
$weekSheet.Cells.Item(419, 4).Value = 526299
$weekSheet.Cells.Item(419, 5).Value = -2.15
$weekSheet.Cells.Item(419, 6).Value = 3032.25
$weekSheet.Cells.Item(419, 7).Value = 371508
$weekSheet.Cells.Item(419, 8).Value = "week"
$weekSheet.Cells.Item(419, 9).Value = "25/10/2024 11:35:48"

# Row 420
$weekSheet.Cells.Item(420, 1).Value = 7
$weekSheet.Cells.Item(420, 2).Value = "ASIANPAINT"
$weekSheet.Cells.Item(420, 3).Value = "Asian Paints Limited"
$weekSheet.Cells.Item(420, 4).Value = 500820
$weekSheet.Cells.Item(420, 5).Value = -0.47
$weekSheet.Cells.Item(420, 6).Value = 2958
$weekSheet.Cells.Item(420, 7).Value = 1047769
$weekSheet.Cells.Item(420, 8).Value = "week"
$weekSheet.Cells.Item(420, 9).Value = "25/10/2024 11:35:48"

# Row 421
$weekSheet.Cells.Item(421, 1).Value = 8
$weekSheet.Cells.Item(421, 2).Value = "GODREJPROP"
$weekSheet.Cells.Item(421, 3).Value = "Godrej Properties Limited"
$weekSheet.Cells.Item(421, 4).Value = 533150
$weekSheet.Cells.Item(421, 5).Value = 0.17
$weekSheet.Cells.Item(421, 6).Value = 2929.25
$weekSheet.Cells.Item(421, 7).Value = 514261
$weekSheet.Cells.Item(421, 8).Value = "week"
$weekSheet.Cells.Item(421, 9).Value = "25/10/2024 11:35:48"

# Row 422
$weekSheet.Cells.Item(422, 1).Value = 9
$weekSheet.Cells.Item(422, 2).Value = "VOLTAS"
$weekSheet.Cells.Item(422, 3).Value = "Voltas Limited"
$weekSheet.Cells.Item(422, 4).Value = 500575
$weekSheet.Cells.Item(422, 5).Value = -2.24
$weekSheet.Cells.Item(422, 6).Value = 1754.85
$weekSheet.Cells.Item(422, 7).Value = 1313221
$weekSheet.Cells.Item(422, 8).Value = "week"
$weekSheet.Cells.Item(422, 9).Value = "25/10/2024 11:35:48"

# Row 423
$weekSheet.Cells.Item(423, 1).Value = 10
$weekSheet.Cells.Item(423, 2).Value = "ICICIBANK"
$weekSheet.Cells.Item(423, 3).Value = "Icici Bank Limited"
$weekSheet.Cells.Item(423, 4).Value = 532174
$weekSheet.Cells.Item(423, 5).Value = 0.22
$weekSheet.Cells.Item(423, 6).Value = 1255.45
$weekSheet.Cells.Item(423, 7).Value = 13546055
$weekSheet.Cells.Item(423, 8).Value = "week"
$weekSheet.Cells.Item(423, 9).Value = "25/10/2024 11:35:48"

# Row 424
$weekSheet.Cells.Item(424, 1).Value = 11
$weekSheet.Cells.Item(424, 2).Value = "TATACHEM"
$weekSheet.Cells.Item(424, 3).Value = "Tata Chemicals Limited"
$weekSheet.Cells.Item(424, 4).Value = 500770
$weekSheet.Cells.Item(424, 5).Value = -3.84
$weekSheet.Cells.Item(424, 6).Value = 1064.75
$weekSheet.Cells.Item(424, 7).Value = 1915919
$weekSheet.Cells.Item(424, 8).Value = "week"
$weekSheet.Cells.Item(424, 9).Value = "25/10/2024 11:35:48"

# Row 425
$weekSheet.Cells.Item(425, 1).Value = 12
$weekSheet.Cells.Item(425, 2).Value = "PEL"
$weekSheet.Cells.Item(425, 3).Value = "Piramal Enterprises Limited"
$weekSheet.Cells.Item(425, 4).Value = 500302
$weekSheet.Cells.Item(425, 5).Value = -0.04
$weekSheet.Cells.Item(425, 6).Value = 1051.8
$weekSheet.Cells.Item(425, 7).Value = 2178692
$weekSheet.Cells.Item(425, 8).Value = "week"
$weekSheet.Cells.Item(425, 9).Value = "25/10/2024 11:35:48"

# Row 426
$weekSheet.Cells.Item(426, 1).Value = 13
$weekSheet.Cells.Item(426, 2).Value = "INDUSINDBK"
$weekSheet.Cells.Item(426, 3).Value = "Indusind Bank Limited"
$weekSheet.Cells.Item(426, 4).Value = 532187
$weekSheet.Cells.Item(426, 5).Value = -18.63
$weekSheet.Cells.Item(426, 6).Value = 1041.6
$weekSheet.Cells.Item(426, 7).Value = 56935420
$weekSheet.Cells.Item(426, 8).Value = "week"
$weekSheet.Cells.Item(426, 9).Value = "25/10/2024 11:35:48"

# Row 427
$weekSheet.Cells.Item(427, 1).Value = 14
$weekSheet.Cells.Item(427, 2).Value = "TATACONSUM"
$weekSheet.Cells.Item(427, 3).Value = "TATA Consumer Products Ltd"
$weekSheet.Cells.Item(427, 4).Value = 500800
$weekSheet.Cells.Item(427, 5).Value = -2.35
$weekSheet.Cells.Item(427, 6).Value = 973.05
$weekSheet.Cells.Item(427, 7).Value = 2258366
$weekSheet.Cells.Item(427, 8).Value = "week"
$weekSheet.Cells.Item(427, 9).Value = "25/10/2024 11:35:48"

# Row 428
$weekSheet.Cells.Item(428, 1).Value = 15
$weekSheet.Cells.Item(428, 2).Value = "JINDALSTEL"
$weekSheet.Cells.Item(428, 3).Value = "Jindal Steel & Power Limited"
$weekSheet.Cells.Item(428, 4).Value = 532286
$weekSheet.Cells.Item(428, 5).Value = -2.73
$weekSheet.Cells.Item(428, 6).Value = 899.05
$weekSheet.Cells.Item(428, 7).Value = 2311226
$weekSheet.Cells.Item(428, 8).Value = "week"
$weekSheet.Cells.Item(428, 9).Value = "25/10/2024 11:35:48"

# Row 429
$weekSheet.Cells.Item(429, 1).Value = 16
$weekSheet.Cells.Item(429, 2).Value = "SYNGENE"
$weekSheet.Cells.Item(429, 3).Value = "Syngene International Limited"
$weekSheet.Cells.Item(429, 4).Value = 539268
$weekSheet.Cells.Item(429, 5).Value = -0.56
$weekSheet.Cells.Item(429, 6).Value = 874.85
$weekSheet.Cells.Item(429, 7).Value = 1367217
$weekSheet.Cells.Item(429, 8).Value = "week"
$weekSheet.Cells.Item(429, 9).Value = "25/10/2024 11:35:48"

# Row 430
$weekSheet.Cells.Item(430, 1).Value = 17
$weekSheet.Cells.Item(430, 2).Value = "AUBANK"
$weekSheet.Cells.Item(430, 3).Value = "AU Small Finance Bank"
$weekSheet.Cells.Item(430, 4).Value = 540611
$weekSheet.Cells.Item(430, 5).Value = -6.37
$weekSheet.Cells.Item(430, 6).Value = 604.5
$weekSheet.Cells.Item(430, 7).Value = 6950878
$weekSheet.Cells.Item(430, 8).Value = "week"
$weekSheet.Cells.Item(430, 9).Value = "25/10/2024 11:35:48"

# Row 431
$weekSheet.Cells.Item(431, 1).Value = 18
$weekSheet.Cells.Item(431, 2).Value = "LICHSGFIN"
$weekSheet.Cells.Item(431, 3).Value = "Lic Housing Finance Limited"
$weekSheet.Cells.Item(431, 4).Value = 500253
$weekSheet.Cells.Item(431, 5).Value = -0.33
$weekSheet.Cells.Item(431, 6).Value = 598.25
$weekSheet.Cells.Item(431, 7).Value = 2102682
$weekSheet.Cells.Item(431, 8).Value = "week"
$weekSheet.Cells.Item(431, 9).Value = "25/10/2024 11:35:48"

# Row 432
$weekSheet.Cells.Item(432, 1).Value = 19
$weekSheet.Cells.Item(432, 2).Value = "GRANULES"
$weekSheet.Cells.Item(432, 3).Value = "Granules India Limited"
$weekSheet.Cells.Item(432, 4).Value = 532482
$weekSheet.Cells.Item(432, 5).Value = -3.12
$weekSheet.Cells.Item(432, 6).Value = 539.6
$weekSheet.Cells.Item(432, 7).Value = 2072052
$weekSheet.Cells.Item(432, 8).Value = "week"
$weekSheet.Cells.Item(432, 9).Value = "25/10/2024 11:35:48"

# Row 433
$weekSheet.Cells.Item(433, 1).Value = 20
$weekSheet.Cells.Item(433, 2).Value = "DABUR"
$weekSheet.Cells.Item(433, 3).Value = "Dabur India Limited"
$weekSheet.Cells.Item(433, 4).Value = 500096
$weekSheet.Cells.Item(433, 5).Value = -0.43
$weekSheet.Cells.Item(433, 6).Value = 538.7
$weekSheet.Cells.Item(433, 7).Value = 1507028
$weekSheet.Cells.Item(433, 8).Value = "week"
$weekSheet.Cells.Item(433, 9).Value = "25/10/2024 11:35:48"

# Row 434
$weekSheet.Cells.Item(434, 1).Value = 21
$weekSheet.Cells.Item(434, 2).Value = "BERGEPAINT"
$weekSheet.Cells.Item(434, 3).Value = "Berger Paints (i) Limited"
$weekSheet.Cells.Item(434, 4).Value = 509480
$weekSheet.Cells.Item(434, 5).Value = -1.3
$weekSheet.Cells.Item(434, 6).Value = 536.95
$weekSheet.Cells.Item(434, 7).Value = 1092908
$weekSheet.Cells.Item(434, 8).Value = "week"
$weekSheet.Cells.Item(434, 9).Value = "25/10/2024 11:35:48"

# Row 435
$weekSheet.Cells.Item(435, 1).Value = 22
$weekSheet.Cells.Item(435, 2).Value = "PETRONET"
$weekSheet.Cells.Item(435, 3).Value = "Petronet Lng Limited"
$weekSheet.Cells.Item(435, 4).Value = 532522
$weekSheet.Cells.Item(435, 5).Value = -0.3
$weekSheet.Cells.Item(435, 6).Value = 337.55
$weekSheet.Cells.Item(435, 7).Value = 4261707
$weekSheet.Cells.Item(435, 8).Value = "week"
$weekSheet.Cells.Item(435, 9).Value = "25/10/2024 11:35:48"

# Row 436
$weekSheet.Cells.Item(436, 1).Value = 23
$weekSheet.Cells.Item(436, 2).Value = "BIOCON"
$weekSheet.Cells.Item(436, 3).Value = "Biocon Limited"
$weekSheet.Cells.Item(436, 4).Value = 532523
$weekSheet.Cells.Item(436, 5).Value = -3.3
$weekSheet.Cells.Item(436, 6).Value = 312
$weekSheet.Cells.Item(436, 7).Value = 3282892
$weekSheet.Cells.Item(436, 8).Value = "week"
$weekSheet.Cells.Item(436, 9).Value = "25/10/2024 11:35:48"

Write-Host "Completed updates"